$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: Salario Basico (G16) corrected
$ws.Range("G16").Value = 3500000

# Row 17: Periodo Mora / Valor Mora updated (now matches period 2503)
$ws.Range("E17").Value = "2503"
$ws.Range("F17").Value = 56940

# Row 19: Periodo Mora / Valor Mora updated (now matches period 2505)
$ws.Range("E19").Value = "2505"
$ws.Range("F19").Value = 140000
